$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4), mirroring the structure/format of row 3.

# A4: Date (serial 42635.642928240741) - same date number format as A3
$ws.Range("A4").Value = 42635.642928240741
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"

# B4: Profitable - boolean FALSE
$ws.Range("B4").Value = $false

# C4: Principle
$ws.Range("C4").Value = 9951.5

# D4: Start Principle
$ws.Range("D4").Value = 10000

# E4: BuyPrice
$ws.Range("E4").Value = 309

# F4: SellPrice
$ws.Range("F4").Value = 312

# G4: IsShortSell - boolean TRUE, same number format as G3
$ws.Range("G4").Value = $true
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"

# H4: Price Change %
$ws.Range("H4").Value = 0.97

# I4: Strong trade - boolean FALSE
$ws.Range("I4").Value = $false
